# Update the barycenter documentation notes in every conductor-component
# input sheet (CHAN, STR_MIX, STR_SC, STR_STAB, Z_JACKET): append a remark
# that the X/Y barycenter values are unused when ITYMESH = -1 in the
# conductor_grid file. Also XLENGTH -> ZLENGTH note is part of the same
# commit family (template_conductor_definition, not this workbook) and the
# absolute-path bookkeeping metadata is left to Excel itself.

$wb = $excel.ActiveWorkbook

$suffix = ". Not used if flag ITYMESH = -1 in file conductor_grid."

$sheetNames = @("CHAN", "STR_MIX", "STR_SC", "STR_STAB", "Z_JACKET")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Locate the row that holds the "X_barycenter" label in column A.
    $xRow = 0
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $label = $ws.Cells.Item($r, 1).Text
        if ($label -eq "X_barycenter") {
            $xRow = $r
            break
        }
    }

    $yRow = $xRow + 1

    $xNoteCell = $ws.Cells.Item($xRow, 4)
    $yNoteCell = $ws.Cells.Item($yRow, 4)

    $xNoteCell.Value = "x coordinate of the barycenter" + $suffix
    $yNoteCell.Value = "y coordinate of the barycenter" + $suffix

    # Mimic the natural cursor landing spot after editing the two cells in
    # turn (Enter moves the selection one row down after the last edit).
    $landingRow = $yRow + 1
    $null = $ws.Range("D" + $landingRow).Select()
}

# The STR_STAB sheet was the last one edited, so it ends up active.
$wb.Worksheets.Item("STR_STAB").Activate()
